$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new metric data row (row 25) after the existing data (rows 1-24)
$ws.Range("A25").Value = "2025-04-28 22:08:54"
$ws.Range("B25").Value = 604
